# aggiornamento fino a 20/09/2021
#
# Extend the daily COVID-style time series on Sheet1 with 11 new rows
# (r375:r385), continuing straight on from the last existing row (374,
# date serial 44448 = 2021-09-09) through date serial 44459 (2021-09-20).
# Columns: A = date, B = "nuovi pos.", C = "somma mobile 7gg.",
# D = "somma mobile 7gg. per 100mila abitanti".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44449, 3, 5, 80.11536612722321),
    @(44450, 0, 5, 80.11536612722321),
    @(44451, 0, 5, 80.11536612722321),
    @(44452, 0, 5, 80.11536612722321),
    @(44453, 0, 5, 80.11536612722321),
    @(44454, 0, 5, 80.11536612722321),
    @(44455, 0, 3, 48.06921967633392),
    @(44456, 0, 0, 0),
    @(44457, 0, 0, 0),
    @(44458, 0, 0, 0),
    @(44459, 0, 0, 0)
)

$startRow = 375
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
}

# The date column (A) in the existing data is styled (bold, bordered,
# centered, "YYYY-MM-DD HH:MM:SS" number format) the same way all the way
# down. Re-use that exact formatting for the newly added date cells by
# copying it from the last pre-existing date cell (A374), rather than
# re-deriving each attribute by hand.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Added rows 375-385 (2021-09-10 .. 2021-09-20)"
